$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells: "latitude" -> "lat", "longitude" -> "lon"
$ws.Range("D1").Value = "lat"
$ws.Range("E1").Value = "lon"

# Update the active selection on the sheet view to D2
$ws.Range("D2").Select()
